# fix(excel): set `sheet_name` and `tpl_name` to load the correct sheet
#
# Renames the two worksheets from the generic "sheet0"/"sheet1" placeholders
# to "test"/"another" and adds a matching literal string ("test"/"another")
# into cell A3 of each sheet, just above the existing template rows, so the
# sheet name and the template content agree. Finally the active cell /
# selection on each sheet is moved onto the newly-added row.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "sheet0" -> "test" -------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "test"
$ws1.Range("A3").Value = "test"
$ws1.Range("A3").Select()

# --- Sheet 2: "sheet1" -> "another" -----------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "another"
$ws2.Range("A3").Value = "another"
$ws2.Range("A6").Select()
